$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2358.4287
$ws.Range("M96").Value = -985.4287000000004
$ws.Range("H96").Value = 786.1429000000001
$ws.Range("I96").Value = 786.1429000000001
$ws.Range("L96").Value = 0
$ws.Range("H111").Value = 13344.695
$ws.Range("J111").Value = 3530.5
$ws.Range("L111").Value = 10591.5
$ws.Range("N111").Value = -16725.5
$ws.Range("J136").Value = 44800
$ws.Range("H136").Value = 40055.734
$ws.Range("N136").Value = -55000
$ws.Range("L136").Value = 44800
$ws.Range("L137").Value = 7469.499899999999
$ws.Range("M137").Value = -4920.900000000001
$ws.Range("H137").Value = 2490.1667
$ws.Range("I137").Value = 2490.3
$ws.Range("J137").Value = 2489.8333
$ws.Range("K137").Value = 7470.900000000001
$ws.Range("N137").Value = -12569.4999
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -1446.2778
$ws.Range("I2").Value = 1559.2778
$ws.Range("H2").Value = 1532.4166
$ws.Range("K2").Value = 1559.2778
$ws.Range("H32").Value = 9564.867
$ws.Range("N32").Value = -16405.333
$ws.Range("J32").Value = 15831.333
$ws.Range("L32").Value = 15831.333
$ws.Range("I45").Value = 2594.2778
$ws.Range("J45").Value = 2468.5334
$ws.Range("M45").Value = -2217.2778
$ws.Range("L45").Value = 2468.5334
$ws.Range("N45").Value = -3222.5334
$ws.Range("H45").Value = 2537.121
$ws.Range("K45").Value = 2594.2778
$ws.Range("M61").Value = -6716.36
$ws.Range("K61").Value = 6928.36
$ws.Range("H61").Value = 11709.531
$ws.Range("I61").Value = 6928.36
$ws.Range("I102").Value = 2644.36
$ws.Range("H102").Value = 7409856
$ws.Range("M102").Value = -1022.36
$ws.Range("K102").Value = 2644.36
$ws.Range("I110").Value = 651.75
$ws.Range("K110").Value = 651.75
$ws.Range("M110").Value = 1393.25
$ws.Range("H110").Value = 731
$ws.Range("I116").Value = 1559.2778
$ws.Range("H116").Value = 1532.4166
$ws.Range("M116").Value = 734.7221999999999
$ws.Range("K116").Value = 1559.2778
$ws.Range("K122").Value = 8225.76
$ws.Range("L122").Value = 13232.625
$ws.Range("H122").Value = 3146.5151
$ws.Range("N122").Value = -18132.625
$ws.Range("J122").Value = 4410.875
$ws.Range("M122").Value = -5775.76
$ws.Range("I122").Value = 2741.92
$ws.Range("L132").Value = 6168
$ws.Range("J132").Value = 2056
$ws.Range("H132").Value = 1584.6721
$ws.Range("N132").Value = -11228
$ws.Range("H136").Value = 11709.531
$ws.Range("K136").Value = 20785.08
$ws.Range("I136").Value = 6928.36
$ws.Range("M136").Value = -18235.08

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1532.4166
$ws.Range("M3").Value = -1445.2778
$ws.Range("I3").Value = 1559.2778
$ws.Range("K3").Value = 1559.2778
$ws.Range("H107").Value = 1258.1177
$ws.Range("M107").Value = 807.4000000000001
$ws.Range("I107").Value = 1112.6
$ws.Range("K107").Value = 1112.6
$ws.Range("I134").Value = 4826.857
$ws.Range("M134").Value = -11945.571
$ws.Range("K134").Value = 14480.571
$ws.Range("H134").Value = 5947.647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -2240.9167
$ws.Range("H31").Value = 3143.8823
$ws.Range("K31").Value = 2535.9167
$ws.Range("I31").Value = 2535.9167
$ws.Range("H34").Value = 3143.8823
$ws.Range("M34").Value = -2333.9167
$ws.Range("I34").Value = 2535.9167
$ws.Range("K34").Value = 2535.9167
$ws.Range("J58").Value = 3123.875
$ws.Range("L58").Value = 3123.875
$ws.Range("N58").Value = -3529.875
$ws.Range("H58").Value = 2921.25
$ws.Range("M132").Value = -6670874.600000001
$ws.Range("K132").Value = 6673404.600000001
$ws.Range("I132").Value = 2224468.2
$ws.Range("H132").Value = 1293199.2
$ws.Range("J136").Value = 3123.875
$ws.Range("H136").Value = 2921.25
$ws.Range("N136").Value = -14471.625
$ws.Range("L136").Value = 9371.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I102").Value = 6524.9375
$ws.Range("L102").Value = 70470
$ws.Range("N102").Value = -75338
$ws.Range("H102").Value = 11688.218
$ws.Range("J102").Value = 23490
$ws.Range("M102").Value = -17140.8125
$ws.Range("K102").Value = 19574.8125
$ws.Range("M131").Value = -1761809.7
$ws.Range("K131").Value = 1766849.7
$ws.Range("H131").Value = 148593.42
$ws.Range("I131").Value = 588949.9
$ws.Range("L137").Value = 9366
$ws.Range("M137").Value = -3769.5
$ws.Range("H137").Value = 3055.8
$ws.Range("I137").Value = 2956.5
$ws.Range("J137").Value = 3122
$ws.Range("K137").Value = 8869.5
$ws.Range("N137").Value = -19566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J70").Value = 7252
$ws.Range("L70").Value = 7252
$ws.Range("H70").Value = 6964.5713
$ws.Range("N70").Value = -7792
$ws.Range("J73").Value = 7252
$ws.Range("H73").Value = 6964.5713
$ws.Range("L73").Value = 7252
$ws.Range("N73").Value = -9124
$ws.Range("N80").Value = -5247.1538
$ws.Range("L80").Value = 3251.1538
$ws.Range("J80").Value = 3251.1538
$ws.Range("H80").Value = 2969.6
$ws.Range("H83").Value = 2969.6
$ws.Range("N83").Value = -26239.769
$ws.Range("J83").Value = 3251.1538
$ws.Range("L83").Value = 16255.769
$ws.Range("I97").Value = 952.5454999999999
$ws.Range("K97").Value = 952.5454999999999
$ws.Range("M97").Value = -456.5454999999999
$ws.Range("H97").Value = 2051.8948
$ws.Range("H113").Value = 1088
$ws.Range("K113").Value = 1088
$ws.Range("I113").Value = 1088
$ws.Range("M113").Value = 1082
$ws.Range("M132").Value = -20410914.5
$ws.Range("K132").Value = 20413444.5
$ws.Range("I132").Value = 6804481.5
$ws.Range("H132").Value = 6174758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 733.5
$ws.Range("I16").Value = 469.09525
$ws.Range("M16").Value = -299.09525
$ws.Range("K16").Value = 469.09525
$ws.Range("N16").Value = -2184
$ws.Range("J16").Value = 1844
$ws.Range("L16").Value = 1844
$ws.Range("M132").Value = -7933.25
$ws.Range("K132").Value = 10463.25
$ws.Range("I132").Value = 3487.75
$ws.Range("H132").Value = 3560.7122
$ws.Range("H136").Value = 4860.7744
$ws.Range("K136").Value = 12124.4349
$ws.Range("I136").Value = 4041.4783
$ws.Range("M136").Value = -9574.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K100").Value = 3169.2
$ws.Range("H100").Value = 1731.8422
$ws.Range("I100").Value = 1584.6
$ws.Range("M100").Value = -2628.2
$ws.Range("H113").Value = 4906624
$ws.Range("N113").Value = -6544.7
$ws.Range("J113").Value = 734.9
$ws.Range("L113").Value = 2204.7
$ws.Range("J136").Value = 4759.2856
$ws.Range("H136").Value = 6855.829
$ws.Range("K136").Value = 21862.4121
$ws.Range("I136").Value = 7287.4707
$ws.Range("N136").Value = -19377.8568
$ws.Range("M136").Value = -19312.4121
$ws.Range("L136").Value = 14277.8568
